# Generate Report for Archive
# Updates the Status of two files (00ec7677-27a6-427d-8c89-36ba41b0dd4a.md and
# 2967e584-08c7-4bff-850d-71fa8b2513a3.md) from "Ready for handoff" to
# "In Translation" across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: columns B (zh-cn) and C (de-de), rows 4 and 5 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B4").Value = $newStatus
$wsOverview.Range("C4").Value = $newStatus
$wsOverview.Range("B5").Value = $newStatus
$wsOverview.Range("C5").Value = $newStatus

# --- zh-cn sheet: column B (Status), rows 4 and 5 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B4").Value = $newStatus
$wsZhCn.Range("B5").Value = $newStatus

# --- de-de sheet: column B (Status), rows 4 and 5 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B4").Value = $newStatus
$wsDeDe.Range("B5").Value = $newStatus
